{"js": "// Insert a new bulleted list item right after the paragraph that reads\n// \"Salvar arquivos PDF em pasta privada, como hist\u00f3rico.\" in the\n// \"REQUISITOS APRESENTADOS PELO AUTOR\" list, adding the new requirement:\n// \"Campo para inser\u00e7\u00e3o de c\u00f3digo do produto atrav\u00e9s do teclado.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Salvar arquivos PDF em pasta privada, como hist\u00f3rico.\";\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === anchorText) {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  // Fallback: the new bullet belongs at the very end of the document body,\n  // right after the last paragraph of the \"REQUISITOS APRESENTADOS PELO\n  // AUTOR\" list.\n  anchor = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// insertParagraph copies the anchor paragraph's style/formatting (pStyle,\n// numPr, spacing, jc, rPr) onto the newly created paragraph, matching the\n// target markup exactly.\nconst newParagraph = anchor.insertParagraph(\n  \"Campo para inser\u00e7\u00e3o de c\u00f3digo do produto atrav\u00e9s do teclado.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item right after the paragraph that reads\n# \"Salvar arquivos PDF em pasta privada, como hist\u00f3rico.\" in the\n# \"REQUISITOS APRESENTADOS PELO AUTOR\" list, adding the new requirement:\n# \"Campo para inser\u00e7\u00e3o de c\u00f3digo do produto atrav\u00e9s do teclado.\"\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Salvar arquivos PDF em pasta privada, como hist\u00f3rico.\"\n\n$anchor = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r`a`n`v\") -eq $anchorText) {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    # Fallback: new bullet belongs at the very end of the document body.\n    $anchor = $d.Paragraphs.Item($d.Paragraphs.Count)\n}\n\n# InsertParagraphAfter duplicates the anchor paragraph's formatting (pStyle,\n# numPr, spacing, jc, rPr) onto the newly created paragraph, matching the\n# target markup exactly - just like pressing Enter at the end of the bullet.\n$anchor.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph.Range.Text = \"Campo para inser\u00e7\u00e3o de c\u00f3digo do produto atrav\u00e9s do teclado.\"\n"}
